$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - F2..F10
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5521
$wsExhibit.Range("F3").Value = 607
$wsExhibit.Range("F4").Value = 12221
$wsExhibit.Range("F6").Value = 616
$wsExhibit.Range("F7").Value = 183
$wsExhibit.Range("F8").Value = 342
$wsExhibit.Range("F9").Value = 1121
$wsExhibit.Range("F10").Value = 106

# Sheet "全部类型" (All types) - F3,F4,F6,F8,F9,F12,F13,F15
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5521
$wsAll.Range("F4").Value = 607
$wsAll.Range("F6").Value = 12221
$wsAll.Range("F8").Value = 616
$wsAll.Range("F9").Value = 183
$wsAll.Range("F12").Value = 342
$wsAll.Range("F13").Value = 1121
$wsAll.Range("F15").Value = 106
